$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.278478666666667
$ws.Range("H2").Value = 6.835436000000001
$ws.Range("I2").Value = 0.006425134583531504
$ws.Range("J2").Value = 0.006425134583531504
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.7065936666666666
$ws.Range("N2").Value = 2.119781
$ws.Range("O2").Value = 0.005187843618793344
$ws.Range("P2").Value = 0.005187843618793344
$ws.Range("Q2").Value = 1.609958595501778
$ws.Range("R2").Value = 14.489627359516
$ws.Range("S2").Value = 0.00003333259344906234
$ws.Range("T2").Value = 0.00003333259344906234

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.278478666666667
$ws.Range("H3").Value = 6.835436000000001
$ws.Range("I3").Value = 0.006425134583531504
$ws.Range("J3").Value = 0.006425134583531504
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 111.9320066666667
$ws.Range("N3").Value = 335.79602
$ws.Range("O3").Value = 0.8218100075305903
$ws.Range("P3").Value = 0.8218100075305903
$ws.Range("Q3").Value = 255.0346893071911
$ws.Range("R3").Value = 2295.31220376472
$ws.Range("S3").Value = 0.005280239900477082
$ws.Range("T3").Value = 0.005280239900477082

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.278478666666667
$ws.Range("H4").Value = 6.835436000000001
$ws.Range("I4").Value = 0.006425134583531504
$ws.Range("J4").Value = 0.006425134583531504
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 23.563205
$ws.Range("N4").Value = 70.689615
$ws.Range("O4").Value = 0.1730021488506163
$ws.Range("P4").Value = 0.1730021488506163
$ws.Range("Q4").Value = 53.68825991079334
$ws.Range("R4").Value = 483.19433919714
$ws.Range("S4").Value = 0.00111156208960536
$ws.Range("T4").Value = 0.00111156208960536

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 346.0613606666668
$ws.Range("H5").Value = 1038.184082
$ws.Range("I5").Value = 0.975866418664458
$ws.Range("J5").Value = 0.9758664186644579
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.7065936666666666
$ws.Range("N5").Value = 2.119781
$ws.Range("O5").Value = 0.005187843618793344
$ws.Range("P5").Value = 0.005187843618793344
$ws.Range("Q5").Value = 244.5247657251158
$ws.Range("R5").Value = 2200.722891526042
$ws.Range("S5").Value = 0.005062642372863122
$ws.Range("T5").Value = 0.005062642372863121

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 346.0613606666668
$ws.Range("H6").Value = 1038.184082
$ws.Range("I6").Value = 0.975866418664458
$ws.Range("J6").Value = 0.9758664186644579
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 111.9320066666667
$ws.Range("N6").Value = 335.79602
$ws.Range("O6").Value = 0.8218100075305903
$ws.Range("P6").Value = 0.8218100075305903
$ws.Range("Q6").Value = 38735.34252921708
$ws.Range("R6").Value = 348618.0827629537
$ws.Range("S6").Value = 0.8019767888714884
$ws.Range("T6").Value = 0.8019767888714884

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 346.0613606666668
$ws.Range("H7").Value = 1038.184082
$ws.Range("I7").Value = 0.975866418664458
$ws.Range("J7").Value = 0.9758664186644579
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 23.563205
$ws.Range("N7").Value = 70.689615
$ws.Range("O7").Value = 0.1730021488506163
$ws.Range("P7").Value = 0.1730021488506163
$ws.Range("Q7").Value = 8154.314783967606
$ws.Range("R7").Value = 73388.83305570846
$ws.Range("S7").Value = 0.1688269874201064
$ws.Range("T7").Value = 0.1688269874201064

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 6.279762333333333
$ws.Range("H8").Value = 18.839287
$ws.Range("I8").Value = 0.01770844675201047
$ws.Range("J8").Value = 0.01770844675201047
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.7065936666666666
$ws.Range("N8").Value = 2.119781
$ws.Range("O8").Value = 0.005187843618793344
$ws.Range("P8").Value = 0.005187843618793344
$ws.Range("Q8").Value = 4.437240292905222
$ws.Range("R8").Value = 39.93516263614699
$ws.Range("S8").Value = 0.00009186865248115925
$ws.Range("T8").Value = 0.00009186865248115925

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 6.279762333333333
$ws.Range("H9").Value = 18.839287
$ws.Range("I9").Value = 0.01770844675201047
$ws.Range("J9").Value = 0.01770844675201047
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 111.9320066666667
$ws.Range("N9").Value = 335.79602
$ws.Range("O9").Value = 0.8218100075305903
$ws.Range("P9").Value = 0.8218100075305903
$ws.Range("Q9").Value = 702.9063993597489
$ws.Range("R9").Value = 6326.157594237739
$ws.Range("S9").Value = 0.01455297875862479
$ws.Range("T9").Value = 0.01455297875862479

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 6.279762333333333
$ws.Range("H10").Value = 18.839287
$ws.Range("I10").Value = 0.01770844675201047
$ws.Range("J10").Value = 0.01770844675201047
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 23.563205
$ws.Range("N10").Value = 70.689615
$ws.Range("O10").Value = 0.1730021488506163
$ws.Range("P10").Value = 0.1730021488506163
$ws.Range("Q10").Value = 147.9713272116117
$ws.Range("R10").Value = 1331.741944904505
$ws.Range("S10").Value = 0.003063599340904529
$ws.Range("T10").Value = 0.003063599340904529
